$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (A..Q):
# A userid, B courseid, C category, D prenom, E nom, F email, G institution,
# H course full name, I course short name, J code cpne, K facetofaceid,
# L facetofacename, M facetofacecapacite, N facetofacesessionid,
# O datedebutsession, P datefinsession, Q duration_in_hours

$rows = @(
    @{ Row = 2;  A = 11691; D = "Romain";       E = "COUPPE";     F = "Romain.Couppe@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" },
    @{ Row = 3;  A = 11628; D = "Michel";        E = "LAUTRIDOU";  F = "Michel.Lautridou@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" },
    @{ Row = 4;  A = 11786; D = "Stéphane";      E = "VILLETTE";   F = "Stephane.VILLETTE@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" },
    @{ Row = 5;  A = 11929; D = "Hervé";         E = "GUION";      F = "Herve.GUION@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" },
    @{ Row = 6;  A = 11712; D = "Guillaume";     E = "TREBUTIEN";  F = "Guillaume.TREBUTIEN@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" },
    @{ Row = 7;  A = 15809; D = "Jean-Jacques";  E = "MUGABE";     F = "Jean-Jacques.MUGABE@fr.toyota-industries.eu";
       G = "Toyota Material Handling France S.A.S."; N = 2012; O = "2025-09-03 13:30"; P = "2025-09-04 17:30" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A          # A userid
    $ws.Cells.Item($row, 2).Value = 1280           # B courseid
    $ws.Cells.Item($row, 3).Value = 300            # C category
    $ws.Cells.Item($row, 4).Value = $r.D           # D prenom
    $ws.Cells.Item($row, 5).Value = $r.E           # E nom
    $ws.Cells.Item($row, 6).Value = $r.F           # F email
    $ws.Cells.Item($row, 7).Value = $r.G           # G institution
    $ws.Cells.Item($row, 8).Value = "LITHIUM-ION TMHMS & TMHMI"   # H course full name
    $ws.Cells.Item($row, 9).Value = "LITHIUM-ION TMHMS & TMHMI"   # I course short name
    $ws.Cells.Item($row, 10).Value = "1404-T2-TE-61"              # J code cpne
    $ws.Cells.Item($row, 11).Value = 719                          # K facetofaceid
    $ws.Cells.Item($row, 12).Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"  # L facetofacename
    $ws.Cells.Item($row, 13).Value = 6                            # M facetofacecapacite
    $ws.Cells.Item($row, 14).Value = $r.N                         # N facetofacesessionid
    $ws.Cells.Item($row, 15).Value = $r.O                         # O datedebutsession
    $ws.Cells.Item($row, 16).Value = $r.P                         # P datefinsession
    $ws.Cells.Item($row, 17).Value = 12                           # Q duration_in_hours
}
